# Add a new "constant_names" worksheet at the end of the workbook, make it
# the active sheet, and populate it with the HL / H2L constant names.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "constant_names"

$newSheet.Range("A1").Value = "HL"
$newSheet.Range("B1").Value = "H2L"

# Match the authored selection/active cell on the new sheet.
$newSheet.Range("A2").Select() | Out-Null
